$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string behind B10/C10 (si index 18) is edited in place:
# it used to hold the long "Objetivos" paragraph and now holds the docente line.
$ws.Range("B10").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("C10").Value = "5817372 - Simone de Fátima Medeiros Sampaio"

# Drop old rows 13-24 entirely; rebuilt below in the new layout/order (now ending at row 23).
$ws.Rows("13:24").Delete()

# Row 13
$ws.Range("A10:C10").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Rows(13).RowHeight = 60
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14
$ws.Range("A10:C10").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Rows(14).RowHeight = 60
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1)Distillation;2)Absorption;3)Liquid-liquid extraction;4)Adsorption."
$ws.Range("C14").Value = "1)Distillation;2)Absorption;3)Liquid-liquid extraction;4)Adsorption."

# Row 15
$ws.Range("A10:C10").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Rows(15).RowHeight = 120
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 16
$ws.Range("A10:C10").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Rows(16).RowHeight = 120
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1)Distillation: liquid-vapor balance; Separation of binary mixtures: flash and continuous separation (continuous rectification); McCabe and Thiele method; Stage efficiency and overall efficiency; Multicomponent distillation - FUG method;2)Absorption and desorption: types of towers; Solubility of gases in liquids; Operations in parallel and countercurrent stage; Mass transfer rates; Countercurrent Multistage operations;3)Liquid-liquid Extraction: Liquid-liquid equilibrium; Extraction in single stage and multistage; Distribution coefficients;4)Adsorption: fundamentals; Operations in single stage and in continuous contact."
$ws.Range("C16").Value = "1)Distillation: liquid-vapor balance; Separation of binary mixtures: flash and continuous separation (continuous rectification); McCabe and Thiele method; Stage efficiency and overall efficiency; Multicomponent distillation - FUG method;2)Absorption and desorption: types of towers; Solubility of gases in liquids; Operations in parallel and countercurrent stage; Mass transfer rates; Countercurrent Multistage operations;3)Liquid-liquid Extraction: Liquid-liquid equilibrium; Extraction in single stage and multistage; Distribution coefficients;4)Adsorption: fundamentals; Operations in single stage and in continuous contact."

# Row 17
$ws.Range("A10").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Avaliação:"

# Row 18
$ws.Range("A10:C10").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Rows(18).RowHeight = 60
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("C18").Value = "5817372 - Simone de Fátima Medeiros Sampaio"

# Row 19
$ws.Range("A10:C10").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Rows(19).RowHeight = 60
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Range("C19").Value = "Aplicação de 2 provas, P1 e P2."

# Row 20
$ws.Range("A10:C10").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Rows(20).RowHeight = 60
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A média do período será MP = (P1+P2)/2. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham frequência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou frequência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham frequência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("C20").Value = "A média do período será MP = (P1+P2)/2. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham frequência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou frequência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham frequência mínima de 70% serão submetidos ao período de recuperação (regimental)."

# Row 21
$ws.Range("A10:C10").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Rows(21).RowHeight = 120
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."
$ws.Range("C21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."

# Row 22
$ws.Range("A10").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Requisitos:"

# Row 23
$ws.Range("B10:C10").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122)
$ws.Rows(23).RowHeight = 30
$ws.Range("B23").Value = "LOQ4054 -  Fenômenos de Transporte III  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4054 -  Fenômenos de Transporte III  (Requisito fraco)`n"

$excel.CutCopyMode = $false